# Update "想去人数" (column F) counts across the 4 sheets of the
# "广州-漫展信息" workbook, per the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 26886
$ws.Range("F4").Value  = 597
$ws.Range("F6").Value  = 622
$ws.Range("F8").Value  = 557
$ws.Range("F13").Value = 50
$ws.Range("F15").Value = 83
$ws.Range("F16").Value = 447
$ws.Range("F18").Value = 1576
$ws.Range("F19").Value = 223
$ws.Range("F20").Value = 59
$ws.Range("F21").Value = 446
$ws.Range("F23").Value = 117

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 4514
$ws.Range("F3").Value  = 240
$ws.Range("F13").Value = 12
$ws.Range("F14").Value = 17

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5118
$ws.Range("F3").Value = 247

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 5118
$ws.Range("F4").Value  = 247
$ws.Range("F5").Value  = 26886
$ws.Range("F6").Value  = 597
$ws.Range("F7").Value  = 4514
$ws.Range("F9").Value  = 240
$ws.Range("F10").Value = 622
$ws.Range("F19").Value = 557
$ws.Range("F25").Value = 50
$ws.Range("F28").Value = 83
$ws.Range("F29").Value = 12
$ws.Range("F30").Value = 17
$ws.Range("F31").Value = 447
$ws.Range("F34").Value = 1576
$ws.Range("F35").Value = 223
$ws.Range("F37").Value = 59
$ws.Range("F38").Value = 446
$ws.Range("F41").Value = 117
